$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: date moves forward one week (44505 -> 44512)
$ws.Range("D7").Value = 44512

# Row 8: date moves (44509 -> 44505), volume updates (200 -> 300)
$ws.Range("D8").Value = 44505
$ws.Range("M8").Value = 300

# Row 9: date moves (44491 -> 44509), and price columns now carry what
# used to be row 8's prices (new week's data replacing the old one)
$ws.Range("D9").Value = 44509
$ws.Range("N9").Value = 19000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 19500
$ws.Range("S9").Value = 2438

# New row 10: the previous (now displaced) week's record is appended,
# carrying the data that row 9 used to hold
$ws.Range("A10").Value = 4
$ws.Range("B10").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C10").Value = "Los Lagos"
$ws.Range("D10").Value = 44491
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100107
$ws.Range("H10").Value = "Otros"
$ws.Range("I10").Value = 100107002
$ws.Range("J10").Value = "Chirimoya"
$ws.Range("K10").Value = "Cultivar IV Región"
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 18000
$ws.Range("O10").Value = 19000
$ws.Range("P10").Value = 18500
$ws.Range("Q10").Value = "$/bandeja 8 kilos"
$ws.Range("R10").Value = "Provincia de Limarí"
$ws.Range("S10").Value = 2312
$ws.Range("T10").Value = 8

# Match the date-formatted style used by the other rows in column D
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
